$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark that sat right after
#    "The highlights are as follows:" (it is being relocated, see step 3).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Reword the "Two-year-olds ..." highlight:
#      "Two-year-olds could already predict upcoming turns at rates
#       greater than chance."
#    becomes
#      "Two-year-olds already predicted upcoming responses at rates
#       greater than chance."
#    Produced as five runs, mirroring the incremental hand edit: drop
#    "could ", add "ed" after "predict", and swap "turns" for "responses ".
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Two-year-olds could already predict upcoming turns*") {
        $targetPara = $candidate
        break
    }
}

$rng = $targetPara.Range.Duplicate
$rng.End = $rng.End - 1   # exclude the paragraph mark
$rng.Text = "Two-year-olds already predict"

$rng.Collapse(0)
$rng.InsertAfter("ed")
$rng.Font.Name = "Times New Roman"
$rng.Font.NameBi = "Times New Roman"

$rng.Collapse(0)
$rng.InsertAfter(" upcoming ")
$rng.Font.Name = "Times New Roman"
$rng.Font.NameBi = "Times New Roman"

$rng.Collapse(0)
$rng.InsertAfter("responses ")
$rng.Font.Name = "Times New Roman"
$rng.Font.NameBi = "Times New Roman"

$rng.Collapse(0)
$rng.InsertAfter("at rates greater than chance.")
$rng.Font.Name = "Times New Roman"
$rng.Font.NameBi = "Times New Roman"

# ---------------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark at the document's final (empty) paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$d.Bookmarks.Add("_GoBack", $lastPara)
